$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the two input values (sig/tau) that drive the rest of the sheet's formulas
$ws.Range("B4").Value = 0.5
$ws.Range("B5").Value = 1.2

# Move the selection to match the author's recorded cursor position
$ws.Range("B4").Select()
